$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.949.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5020"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.29%  "
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07114"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8839"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07563"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.280"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "88.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008361"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.95%  "
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.008.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.019"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.116.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.445"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.849"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.23%  "
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.091"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.638"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.660"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05112"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.038"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.23%  "
$ws.Range("E35").Value = "  -7.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7209"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.26%  "
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.463"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5270"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.508"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.236"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1460"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4587"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.975"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.88%  "
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.02%  "
